$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.956.47'
$ws.Range('E2').Value = '  +1.55%  '
$ws.Range('D3').Value = '3.908.46'
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').Value = "'483.87"
$ws.Range('E5').Value = '  +2.72%  '
$ws.Range('D6').Value = "'144.77"
$ws.Range('E6').Value = '  -0.30%  '
$ws.Range('E7').Value = '  -0.65%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = "'0.726"
$ws.Range('E9').Value = '  -2.17%  '
$ws.Range('E10').Value = '  +1.00%  '
$ws.Range('D11').Value = "'0.0000353"
$ws.Range('E11').Value = '  +4.16%  '
$ws.Range('D12').Value = "'42.47"
$ws.Range('E12').Value = '  -1.77%  '
$ws.Range('D13').Value = "'10.63"
$ws.Range('E13').Value = '  +1.92%  '
$ws.Range('D14').Value = '4.535.49'
$ws.Range('E14').Value = '  +0.34%  '
$ws.Range('D15').Value = "'14.74"
$ws.Range('E15').Value = '  -2.48%  '
$ws.Range('D16').Value = '3.898.17'
$ws.Range('E16').Value = '  -0.42%  '
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('E18').Value = '  -1.38%  '
$ws.Range('E19').Value = '  -2.66%  '
$ws.Range('D20').Value = '68.073.95'
$ws.Range('E20').Value = '  +1.26%  '
$ws.Range('D21').Value = "'445.43"
$ws.Range('E21').Value = '  +3.02%  '
$ws.Range('D22').Value = "'14.63"
$ws.Range('E22').Value = '  -0.90%  '
$ws.Range('E23').Value = '  +0.24%  '
$ws.Range('D24').Value = "'88.98"
$ws.Range('E24').Value = '  +0.74%  '
$ws.Range('D25').Value = "'11.67"
$ws.Range('E25').Value = '  +15.37%  '
$ws.Range('D26').Value = "'10.93"
$ws.Range('E26').Value = '  +12.83%  '
$ws.Range('D27').Value = "'3.59"
$ws.Range('E27').Value = '  +0.78%  '
$ws.Range('D28').Value = "'38.55"
$ws.Range('E28').Value = '  -0.56%  '
$ws.Range('E29').Value = '  +3.17%  '
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').Value = "'13.37"
$ws.Range('E30').Value = '  -2.47%  '
$ws.Range('D31').Value = "'0.130"
$ws.Range('E31').Value = '  -1.10%  '
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').Value = "'686.19"
$ws.Range('E32').Value = '  -6.73%  '
$ws.Range('E33').Value = '  +2.57%  '
$ws.Range('D34').Value = '0.0₃0937'
$ws.Range('E34').Value = '  +25.72%  '
$ws.Range('D35').Value = "'41.34"
$ws.Range('E35').Value = '  -6.55%  '
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').Value = "'58.86"
$ws.Range('E36').Value = '  +1.25%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').Value = "'5.81"
$ws.Range('E37').Value = '  +8.38%  '
$ws.Range('D38').Value = "'0.149"
$ws.Range('E38').Value = '  -5.83%  '
$ws.Range('E39').Value = '  -0.12%  '
$ws.Range('E40').Value = '  -0.30%  '
$ws.Range('D41').Value = "'2.84"
$ws.Range('E41').Value = '  +15.21%  '
$ws.Range('D42').Value = "'3.02"
$ws.Range('E42').Value = '  -5.95%  '
$ws.Range('D43').Value = "'3.02"
$ws.Range('E43').Value = '  +7.57%  '
$ws.Range('D44').Value = "'0.357"
$ws.Range('E44').Value = '  +6.07%  '
$ws.Range('E45').Value = '  +0.30%  '
$ws.Range('E46').Value = '  -0.06%  '
$ws.Range('D47').Value = "'3.43"
$ws.Range('E47').Value = '  -0.03%  '
$ws.Range('D48').Value = "'2.12"
$ws.Range('E48').Value = '  -3.09%  '
$ws.Range('D49').Value = "'146.29"
$ws.Range('E49').Value = '  +2.14%  '
$ws.Range('D50').Value = "'3.15"
$ws.Range('E50').Value = '  -0.77%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0334'
$ws.Range('E51').Value = '  +47.96%  '
